# Applies the NATMI TPM re-calculation update for Cp-Slc40a1 LR-pair sheet.
# Updates columns G:T (expression / specificity / edge-weight metrics) for data rows 2-26
# with the refreshed TPM-derived values, leaving columns A:F (cluster/gene/cell-count) unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 3.323546
$ws.Cells.Item(2, 8).Value = 9.970638000000001
$ws.Cells.Item(2, 9).Value = 0.05491115868684463
$ws.Cells.Item(2, 10).Value = 0.05617460731729474
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.3231746666666667
$ws.Cells.Item(2, 14).Value = 0.9695240000000001
$ws.Cells.Item(2, 15).Value = 0.01551784106606065
$ws.Cells.Item(2, 16).Value = 0.01834699745525213
$ws.Cells.Item(2, 17).Value = 1.074085870701333
$ws.Cells.Item(2, 18).Value = 9.666772836312001
$ws.Cells.Item(2, 19).Value = 0.0008521026332556908
$ws.Cells.Item(2, 20).Value = 0.001030635377500194

# Row 3
$ws.Cells.Item(3, 7).Value = 3.323546
$ws.Cells.Item(3, 8).Value = 9.970638000000001
$ws.Cells.Item(3, 9).Value = 0.05491115868684463
$ws.Cells.Item(3, 10).Value = 0.05617460731729474
$ws.Cells.Item(3, 15).Value = 0.2047034908484057
$ws.Cells.Item(3, 16).Value = 0.2420242873791943
$ws.Cells.Item(3, 17).Value = 14.16879617902467
$ws.Cells.Item(3, 18).Value = 127.519165611222
$ws.Cells.Item(3, 19).Value = 0.01124050586972785
$ws.Cells.Item(3, 20).Value = 0.01359561930477433

# Row 4
$ws.Cells.Item(4, 7).Value = 3.323546
$ws.Cells.Item(4, 8).Value = 9.970638000000001
$ws.Cells.Item(4, 9).Value = 0.05491115868684463
$ws.Cells.Item(4, 10).Value = 0.05617460731729474
$ws.Cells.Item(4, 13).Value = 1.969873333333333
$ws.Cells.Item(4, 14).Value = 5.909619999999999
$ws.Cells.Item(4, 15).Value = 0.09458718290708983
$ws.Cells.Item(4, 16).Value = 0.1118319743518542
$ws.Cells.Item(4, 17).Value = 6.546964637506667
$ws.Cells.Item(4, 18).Value = 58.92268173756
$ws.Cells.Item(4, 19).Value = 0.005193891810352807
$ws.Cells.Item(4, 20).Value = 0.006282117244733186

# Row 5
$ws.Cells.Item(5, 7).Value = 3.323546
$ws.Cells.Item(5, 8).Value = 9.970638000000001
$ws.Cells.Item(5, 9).Value = 0.05491115868684463
$ws.Cells.Item(5, 10).Value = 0.05617460731729474
$ws.Cells.Item(5, 13).Value = 9.634278999999999
$ws.Cells.Item(5, 14).Value = 19.268558
$ws.Cells.Item(5, 15).Value = 0.4626080746059482
$ws.Cells.Item(5, 16).Value = 0.3646327317244112
$ws.Cells.Item(5, 17).Value = 32.019969433334
$ws.Cells.Item(5, 18).Value = 192.119816600004
$ws.Cells.Item(5, 19).Value = 0.02540234539450288
$ws.Cells.Item(5, 20).Value = 0.02048310051965128

# Row 6
$ws.Cells.Item(6, 7).Value = 3.323546
$ws.Cells.Item(6, 8).Value = 9.970638000000001
$ws.Cells.Item(6, 9).Value = 0.05491115868684463
$ws.Cells.Item(6, 10).Value = 0.05617460731729474
$ws.Cells.Item(6, 13).Value = 4.635523666666667
$ws.Cells.Item(6, 14).Value = 13.906571
$ws.Cells.Item(6, 15).Value = 0.2225834105724956
$ws.Cells.Item(6, 16).Value = 0.2631640090892882
$ws.Cells.Item(6, 17).Value = 15.40637614025534
$ws.Cells.Item(6, 18).Value = 138.657385262298
$ws.Cells.Item(6, 19).Value = 0.01222231297900539
$ws.Cells.Item(6, 20).Value = 0.01478313487063575

# Row 7
$ws.Cells.Item(7, 9).Value = 0.6804732481664908
$ws.Cells.Item(7, 10).Value = 0.6961302296255228
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.3231746666666667
$ws.Cells.Item(7, 14).Value = 0.9695240000000001
$ws.Cells.Item(7, 15).Value = 0.01551784106606065
$ws.Cells.Item(7, 16).Value = 0.01834699745525213
$ws.Cells.Item(7, 17).Value = 13.31034927552844
$ws.Cells.Item(7, 18).Value = 119.793143479756
$ws.Cells.Item(7, 19).Value = 0.01055947571475365
$ws.Cells.Item(7, 20).Value = 0.01277189955146355

# Row 8
$ws.Cells.Item(8, 9).Value = 0.6804732481664908
$ws.Cells.Item(8, 10).Value = 0.6961302296255228
$ws.Cells.Item(8, 15).Value = 0.2047034908484057
$ws.Cells.Item(8, 16).Value = 0.2420242873791943
$ws.Cells.Item(8, 19).Value = 0.1392952493286342
$ws.Cells.Item(8, 20).Value = 0.1684804227482321

# Row 9
$ws.Cells.Item(9, 9).Value = 0.6804732481664908
$ws.Cells.Item(9, 10).Value = 0.6961302296255228
$ws.Cells.Item(9, 13).Value = 1.969873333333333
$ws.Cells.Item(9, 14).Value = 5.909619999999999
$ws.Cells.Item(9, 15).Value = 0.09458718290708983
$ws.Cells.Item(9, 16).Value = 0.1118319743518542
$ws.Cells.Item(9, 17).Value = 81.13167521964222
$ws.Cells.Item(9, 18).Value = 730.1850769767799
$ws.Cells.Item(9, 19).Value = 0.06436404758770539
$ws.Cells.Item(9, 20).Value = 0.07784961798503183

# Row 10
$ws.Cells.Item(10, 9).Value = 0.6804732481664908
$ws.Cells.Item(10, 10).Value = 0.6961302296255228
$ws.Cells.Item(10, 13).Value = 9.634278999999999
$ws.Cells.Item(10, 14).Value = 19.268558
$ws.Cells.Item(10, 15).Value = 0.4626080746059482
$ws.Cells.Item(10, 16).Value = 0.3646327317244112
$ws.Cells.Item(10, 17).Value = 396.7997239095336
$ws.Cells.Item(10, 18).Value = 2380.798343457202
$ws.Cells.Item(10, 19).Value = 0.3147924191551559
$ws.Cells.Item(10, 20).Value = 0.253831867264296

# Row 11
$ws.Cells.Item(11, 9).Value = 0.6804732481664908
$ws.Cells.Item(11, 10).Value = 0.6961302296255228
$ws.Cells.Item(11, 13).Value = 4.635523666666667
$ws.Cells.Item(11, 14).Value = 13.906571
$ws.Cells.Item(11, 15).Value = 0.2225834105724956
$ws.Cells.Item(11, 16).Value = 0.2631640090892882
$ws.Cells.Item(11, 17).Value = 190.9197887158388
$ws.Cells.Item(11, 18).Value = 1718.278098442549
$ws.Cells.Item(11, 19).Value = 0.1514620563802417
$ws.Cells.Item(11, 20).Value = 0.1831964220764994

# Row 12
$ws.Cells.Item(12, 7).Value = 4.901883
$ws.Cells.Item(12, 8).Value = 14.705649
$ws.Cells.Item(12, 9).Value = 0.08098822019534135
$ws.Cells.Item(12, 10).Value = 0.0828516748798791
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 0.3231746666666667
$ws.Cells.Item(12, 14).Value = 0.9695240000000001
$ws.Cells.Item(12, 15).Value = 0.01551784106606065
$ws.Cells.Item(12, 16).Value = 0.01834699745525213
$ws.Cells.Item(12, 17).Value = 1.584164404564
$ws.Cells.Item(12, 18).Value = 14.257479641076
$ws.Cells.Item(12, 19).Value = 0.001256762329214431
$ws.Cells.Item(12, 20).Value = 0.001520079468184518

# Row 13
$ws.Cells.Item(13, 7).Value = 4.901883
$ws.Cells.Item(13, 8).Value = 14.705649
$ws.Cells.Item(13, 9).Value = 0.08098822019534135
$ws.Cells.Item(13, 10).Value = 0.0828516748798791
$ws.Cells.Item(13, 15).Value = 0.2047034908484057
$ws.Cells.Item(13, 16).Value = 0.2420242873791943
$ws.Cells.Item(13, 17).Value = 20.897493556709
$ws.Cells.Item(13, 18).Value = 188.077442010381
$ws.Cells.Item(13, 19).Value = 0.01657857139158573
$ws.Cells.Item(13, 20).Value = 0.02005211757097543

# Row 14
$ws.Cells.Item(14, 7).Value = 4.901883
$ws.Cells.Item(14, 8).Value = 14.705649
$ws.Cells.Item(14, 9).Value = 0.08098822019534135
$ws.Cells.Item(14, 10).Value = 0.0828516748798791
$ws.Cells.Item(14, 13).Value = 1.969873333333333
$ws.Cells.Item(14, 14).Value = 5.909619999999999
$ws.Cells.Item(14, 15).Value = 0.09458718290708983
$ws.Cells.Item(14, 16).Value = 0.1118319743518542
$ws.Cells.Item(14, 17).Value = 9.656088604819999
$ws.Cells.Item(14, 18).Value = 86.90479744337999
$ws.Cells.Item(14, 19).Value = 0.007660447596936419
$ws.Cells.Item(14, 20).Value = 0.009265466380174801

# Row 15
$ws.Cells.Item(15, 7).Value = 4.901883
$ws.Cells.Item(15, 8).Value = 14.705649
$ws.Cells.Item(15, 9).Value = 0.08098822019534135
$ws.Cells.Item(15, 10).Value = 0.0828516748798791
$ws.Cells.Item(15, 13).Value = 9.634278999999999
$ws.Cells.Item(15, 14).Value = 19.268558
$ws.Cells.Item(15, 15).Value = 0.4626080746059482
$ws.Cells.Item(15, 16).Value = 0.3646327317244112
$ws.Cells.Item(15, 17).Value = 47.22610844735699
$ws.Cells.Item(15, 18).Value = 283.356650684142
$ws.Cells.Item(15, 19).Value = 0.03746580461032943
$ws.Cells.Item(15, 20).Value = 0.0302104325393931

# Row 16
$ws.Cells.Item(16, 7).Value = 4.901883
$ws.Cells.Item(16, 8).Value = 14.705649
$ws.Cells.Item(16, 9).Value = 0.08098822019534135
$ws.Cells.Item(16, 10).Value = 0.0828516748798791
$ws.Cells.Item(16, 13).Value = 4.635523666666667
$ws.Cells.Item(16, 14).Value = 13.906571
$ws.Cells.Item(16, 15).Value = 0.2225834105724956
$ws.Cells.Item(16, 16).Value = 0.2631640090892882
$ws.Cells.Item(16, 17).Value = 22.722794657731
$ws.Cells.Item(16, 18).Value = 204.505151919579
$ws.Cells.Item(16, 19).Value = 0.01802663426727534
$ws.Cells.Item(16, 20).Value = 0.02180357892115125

# Row 17
$ws.Cells.Item(17, 7).Value = 4.083945
$ws.Cells.Item(17, 8).Value = 8.16789
$ws.Cells.Item(17, 9).Value = 0.06747436381604036
$ws.Cells.Item(17, 10).Value = 0.04601791915029495
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.3231746666666667
$ws.Cells.Item(17, 14).Value = 0.9695240000000001
$ws.Cells.Item(17, 15).Value = 0.01551784106606065
$ws.Cells.Item(17, 16).Value = 0.01834699745525213
$ws.Cells.Item(17, 17).Value = 1.31982756406
$ws.Cells.Item(17, 18).Value = 7.918965384360001
$ws.Cells.Item(17, 19).Value = 0.001047056453730868
$ws.Cells.Item(17, 20).Value = 0.0008442906455464595

# Row 18
$ws.Cells.Item(18, 7).Value = 4.083945
$ws.Cells.Item(18, 8).Value = 8.16789
$ws.Cells.Item(18, 9).Value = 0.06747436381604036
$ws.Cells.Item(18, 10).Value = 0.04601791915029495
$ws.Cells.Item(18, 15).Value = 0.2047034908484057
$ws.Cells.Item(18, 16).Value = 0.2420242873791943
$ws.Cells.Item(18, 17).Value = 17.410495991735
$ws.Cells.Item(18, 18).Value = 104.46297595041
$ws.Cells.Item(18, 19).Value = 0.01381223781591881
$ws.Cells.Item(18, 20).Value = 0.01113745408902351

# Row 19
$ws.Cells.Item(19, 7).Value = 4.083945
$ws.Cells.Item(19, 8).Value = 8.16789
$ws.Cells.Item(19, 9).Value = 0.06747436381604036
$ws.Cells.Item(19, 10).Value = 0.04601791915029495
$ws.Cells.Item(19, 13).Value = 1.969873333333333
$ws.Cells.Item(19, 14).Value = 5.909619999999999
$ws.Cells.Item(19, 15).Value = 0.09458718290708983
$ws.Cells.Item(19, 16).Value = 0.1118319743518542
$ws.Cells.Item(19, 17).Value = 8.0448543503
$ws.Cells.Item(19, 18).Value = 48.2691261018
$ws.Cells.Item(19, 19).Value = 0.006382209991807332
$ws.Cells.Item(19, 20).Value = 0.005146274754141484

# Row 20
$ws.Cells.Item(20, 7).Value = 4.083945
$ws.Cells.Item(20, 8).Value = 8.16789
$ws.Cells.Item(20, 9).Value = 0.06747436381604036
$ws.Cells.Item(20, 10).Value = 0.04601791915029495
$ws.Cells.Item(20, 13).Value = 9.634278999999999
$ws.Cells.Item(20, 14).Value = 19.268558
$ws.Cells.Item(20, 15).Value = 0.4626080746059482
$ws.Cells.Item(20, 16).Value = 0.3646327317244112
$ws.Cells.Item(20, 17).Value = 39.345865550655
$ws.Cells.Item(20, 18).Value = 157.38346220262
$ws.Cells.Item(20, 19).Value = 0.03121418553019968
$ws.Cells.Item(20, 20).Value = 0.01677963956804514

# Row 21
$ws.Cells.Item(21, 7).Value = 4.083945
$ws.Cells.Item(21, 8).Value = 8.16789
$ws.Cells.Item(21, 9).Value = 0.06747436381604036
$ws.Cells.Item(21, 10).Value = 0.04601791915029495
$ws.Cells.Item(21, 13).Value = 4.635523666666667
$ws.Cells.Item(21, 14).Value = 13.906571
$ws.Cells.Item(21, 15).Value = 0.2225834105724956
$ws.Cells.Item(21, 16).Value = 0.2631640090892882
$ws.Cells.Item(21, 17).Value = 18.931223700865
$ws.Cells.Item(21, 18).Value = 113.58734220519
$ws.Cells.Item(21, 19).Value = 0.01501867402438365
$ws.Cells.Item(21, 20).Value = 0.01211026009353835

# Row 22
$ws.Cells.Item(22, 7).Value = 7.030262666666666
$ws.Cells.Item(22, 8).Value = 21.090788
$ws.Cells.Item(22, 9).Value = 0.1161530091352829
$ws.Cells.Item(22, 10).Value = 0.1188255690270083
$ws.Cells.Item(22, 11).Value = 3
$ws.Cells.Item(22, 12).Value = 1
$ws.Cells.Item(22, 13).Value = 0.3231746666666667
$ws.Cells.Item(22, 14).Value = 0.9695240000000001
$ws.Cells.Item(22, 15).Value = 0.01551784106606065
$ws.Cells.Item(22, 16).Value = 0.01834699745525213
$ws.Cells.Item(22, 17).Value = 2.272002793879111
$ws.Cells.Item(22, 18).Value = 20.448025144912
$ws.Cells.Item(22, 19).Value = 0.001802443935106011
$ws.Cells.Item(22, 20).Value = 0.002180092412557407

# Row 23
$ws.Cells.Item(23, 7).Value = 7.030262666666666
$ws.Cells.Item(23, 8).Value = 21.090788
$ws.Cells.Item(23, 9).Value = 0.1161530091352829
$ws.Cells.Item(23, 10).Value = 0.1188255690270083
$ws.Cells.Item(23, 15).Value = 0.2047034908484057
$ws.Cells.Item(23, 16).Value = 0.2420242873791943
$ws.Cells.Item(23, 17).Value = 29.97110881239689
$ws.Cells.Item(23, 18).Value = 269.739979311572
$ws.Cells.Item(23, 19).Value = 0.02377692644253916
$ws.Cells.Item(23, 20).Value = 0.02875867366618895

# Row 24
$ws.Cells.Item(24, 7).Value = 7.030262666666666
$ws.Cells.Item(24, 8).Value = 21.090788
$ws.Cells.Item(24, 9).Value = 0.1161530091352829
$ws.Cells.Item(24, 10).Value = 0.1188255690270083
$ws.Cells.Item(24, 13).Value = 1.969873333333333
$ws.Cells.Item(24, 14).Value = 5.909619999999999
$ws.Cells.Item(24, 15).Value = 0.09458718290708983
$ws.Cells.Item(24, 16).Value = 0.1118319743518542
$ws.Cells.Item(24, 17).Value = 13.84872695339555
$ws.Cells.Item(24, 18).Value = 124.63854258056
$ws.Cells.Item(24, 19).Value = 0.01098658592028788
$ws.Cells.Item(24, 20).Value = 0.01328849798777287

# Row 25
$ws.Cells.Item(25, 7).Value = 7.030262666666666
$ws.Cells.Item(25, 8).Value = 21.090788
$ws.Cells.Item(25, 9).Value = 0.1161530091352829
$ws.Cells.Item(25, 10).Value = 0.1188255690270083
$ws.Cells.Item(25, 13).Value = 9.634278999999999
$ws.Cells.Item(25, 14).Value = 19.268558
$ws.Cells.Item(25, 15).Value = 0.4626080746059482
$ws.Cells.Item(25, 16).Value = 0.3646327317244112
$ws.Cells.Item(25, 17).Value = 67.73151197395066
$ws.Cells.Item(25, 18).Value = 406.389071843704
$ws.Cells.Item(25, 19).Value = 0.05373331991576031
$ws.Cells.Item(25, 20).Value = 0.04332769183302562

# Row 26
$ws.Cells.Item(26, 7).Value = 7.030262666666666
$ws.Cells.Item(26, 8).Value = 21.090788
$ws.Cells.Item(26, 9).Value = 0.1161530091352829
$ws.Cells.Item(26, 10).Value = 0.1188255690270083
$ws.Cells.Item(26, 13).Value = 4.635523666666667
$ws.Cells.Item(26, 14).Value = 13.906571
$ws.Cells.Item(26, 15).Value = 0.2225834105724956
$ws.Cells.Item(26, 16).Value = 0.2631640090892882
$ws.Cells.Item(26, 17).Value = 32.58894897421644
$ws.Cells.Item(26, 18).Value = 293.300540767948
$ws.Cells.Item(26, 19).Value = 0.03127061312746345
$ws.Cells.Item(26, 20).Value = 0.03127061312746345

